# Apply "Tarefa 1 e 2" changes to resultado_comissoes workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab
$ws.Name = "Comissão"

# Update header text in E1
$ws.Range("E1").Value = "Comissão Final"

# Update commission values (Tarefa 1/2 adjustments)
$ws.Range("B3").Value = 160
$ws.Range("B6").Value = 320
$ws.Range("B9").Value = 500
$ws.Range("B13").Value = 150
